$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62: Gridx (B62) bumped 16 -> 18 ---
$ws.Cells.Item(62, 2).Value = 18

# --- Rows 64-80: rewrite the ObjectID/Gridx/Gridy/ObjectName table ---
# Columns: A=ObjectID(number), B=Gridx(number), C=Gridy(number), D=ObjectName(string)
$rows = @(
    @{ Row = 64; A = 0; B = 0;  C = 1; D = "Defense_Generator_0" },
    @{ Row = 65; A = 0; B = 0;  C = 5; D = "Defense_Generator_1" },
    @{ Row = 66; A = 5; B = 8;  C = 0; D = "Portal" },
    @{ Row = 67; A = 5; B = 2;  C = 1; D = "Portal" },
    @{ Row = 68; A = 5; B = 14; C = 6; D = "Portal" },
    @{ Row = 69; A = 2; B = 0;  C = 0; D = "Yama" },
    @{ Row = 70; A = 2; B = 0;  C = 1; D = "Yama" },
    @{ Row = 71; A = 2; B = 2;  C = 0; D = "Yama" },
    @{ Row = 72; A = 2; B = 12; C = 0; D = "Yama" },
    @{ Row = 73; A = 2; B = 14; C = 1; D = "Yama" },
    @{ Row = 74; A = 2; B = 14; C = 0; D = "Yama" },
    @{ Row = 75; A = 3; B = 0;  C = 6; D = "Building" },
    @{ Row = 76; A = 3; B = 1;  C = 6; D = "Building" },
    @{ Row = 77; A = 3; B = 2;  C = 6; D = "Building" },
    @{ Row = 78; A = 3; B = 3;  C = 6; D = "Building" },
    @{ Row = 79; A = 3; B = 15; C = 6; D = "Building" },
    @{ Row = 80; A = 3; B = 15; C = 5; D = "Building" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# --- Sheet view / selection tweak (tooltip position fix) ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("H60").Select()
